$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "pnl" header in H1, matching formatting of existing headers ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "pnl"

# --- Apply header-like style (bold/border/center) to column A data cells A2:A9 ---
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)

# --- Apply the date style (used in column B) to the new C4 / C7 close-date cells ---
$ws.Range("B2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)

# Row 2 : ADA/USDT:USDT long  (trade_id 1)
$ws.Range("B2").Value2 = 44730.03644305556
$ws.Range("F2").Value2 = 0.46076

# Row 3 : APE/USDT:USDT short (trade_id 2)
$ws.Range("B3").Value2 = 44730.03646423611
$ws.Range("E3").Value2 = "short"
$ws.Range("F3").Value2 = 3.331

# Row 4 : AVAX/USDT:USDT long, now closed (trade_id 3)
$ws.Range("B4").Value2 = 44730.0364677662
$ws.Range("C4").Value2 = 44730.04080313955
$ws.Range("F4").Value2 = 15.01
$ws.Range("G4").Value2 = 14.96
$ws.Range("H4").Value2 = -0.3331112591605525

# Row 5 : BTC/USDT:USDT long (trade_id 4)
$ws.Range("B5").Value2 = 44730.03647100695
$ws.Range("F5").Value2 = 19285

# Row 6 : ETH/USDT:USDT long (trade_id 5)
$ws.Range("B6").Value2 = 44730.03647425926
$ws.Range("F6").Value2 = 1006.05

# Row 7 : FTM/USDT:USDT long, now closed (trade_id 6)
$ws.Range("B7").Value2 = 44730.0364772338
$ws.Range("C7").Value2 = 44730.03936880787
$ws.Range("F7").Value2 = 0.22
$ws.Range("G7").Value2 = 0.219
$ws.Range("H7").Value2 = -0.454545454545455

# Row 8 : SOL/USDT:USDT long (trade_id 7)
$ws.Range("B8").Value2 = 44730.03648079861
$ws.Range("F8").Value2 = 28.599

# Row 9 : XRP/USDT:USDT long (trade_id 8)
$ws.Range("B9").Value2 = 44730.03649108797
$ws.Range("F9").Value2 = 0.3091

# --- New rows 10 and 11, formatted like the rest of the data rows ---
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A11:G11").PasteSpecial(-4122)

# Row 10 : FTM/USDT:USDT short - closing leg of the FTM long trade (trade_id 9)
$ws.Range("A10").Value2 = 9
$ws.Range("B10").Value2 = 44730.03936880787
$ws.Range("C10").Value2 = ""
$ws.Range("D10").Value2 = "FTM/USDT:USDT"
$ws.Range("E10").Value2 = "short"
$ws.Range("F10").Value2 = 0.219
$ws.Range("G10").Value2 = ""

# Row 11 : AVAX/USDT:USDT short - closing leg of the AVAX long trade (trade_id 10)
$ws.Range("A11").Value2 = 10
$ws.Range("B11").Value2 = 44730.04080313657
$ws.Range("C11").Value2 = ""
$ws.Range("D11").Value2 = "AVAX/USDT:USDT"
$ws.Range("E11").Value2 = "short"
$ws.Range("F11").Value2 = 14.96
$ws.Range("G11").Value2 = ""

Write-Output "done"
